$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1, 2).Value = "Хэш"
$ws.Cells.Item(1, 3).Value = "Время обработки"
$ws.Cells.Item(1, 4).Value = "Хэммингово расстояние"

# Set column B (hash) to Text format so long digit strings are not
# misinterpreted as numbers in scientific notation.
$ws.Range("B2:B63").NumberFormat = "@"

$ws.Cells.Item(2, 2).Value = "1110101000000000110000000000000010000110000000001000000000000000"
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(3, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(4, 2).Value = "1010100000000000100000000000000010000000000000001000000000000000"
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 6
$ws.Cells.Item(5, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(5, 3).Value = 0.015627
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(6, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(7, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(8, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(9, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(9, 3).Value = 0
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(10, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(11, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(12, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(12, 3).Value = 0
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(13, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(13, 3).Value = 0
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(14, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(15, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(16, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(16, 3).Value = 0
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(17, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(17, 3).Value = 0.015622
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(18, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(19, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(20, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(21, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(21, 3).Value = 0
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(22, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(22, 3).Value = 0
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(23, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(23, 3).Value = 0
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(24, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(24, 3).Value = 0
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(25, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(25, 3).Value = 0.015635
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(26, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(26, 3).Value = 0
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(27, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(27, 3).Value = 0
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(28, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(28, 3).Value = 0
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(29, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(29, 3).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(30, 2).Value = "1110101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(31, 2).Value = "1010101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(31, 3).Value = 0.015627
$ws.Cells.Item(31, 4).Value = 1
$ws.Cells.Item(32, 2).Value = "1010101000000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 1
$ws.Cells.Item(33, 2).Value = "1010101100000000110000001000000010000110000000001000000000000000"
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(33, 4).Value = 2
$ws.Cells.Item(34, 2).Value = "1010000000000000100000000000000000000000000000000000000000000000"
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 9
$ws.Cells.Item(35, 2).Value = "1110101101000000101001000001000010100100000000001000011000000000"
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(35, 4).Value = 11
$ws.Cells.Item(36, 2).Value = "1110100100000000100000000000000010000010100000000000111000000000"
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(36, 4).Value = 10
$ws.Cells.Item(37, 2).Value = "1110100000000000100100001000000000000000000000000000000001000000"
$ws.Cells.Item(37, 3).Value = 0.015625
$ws.Cells.Item(37, 4).Value = 8
$ws.Cells.Item(38, 2).Value = "1110110001100000101000001000100110000000000000001000000000000000"
$ws.Cells.Item(38, 3).Value = 0
$ws.Cells.Item(38, 4).Value = 10
$ws.Cells.Item(39, 2).Value = "1011100000000000110000000000000000000000000000000000000000000000"
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = 8
$ws.Cells.Item(40, 2).Value = "1010100100000000100000000000000010000010000000001000000000000000"
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 6
$ws.Cells.Item(41, 2).Value = "1011100000000000110000001000100000000000100000000000000000000000"
$ws.Cells.Item(41, 3).Value = 0.015626
$ws.Cells.Item(41, 4).Value = 9
$ws.Cells.Item(42, 2).Value = "1011000000000000100000000000000000000000000000000000000000000000"
$ws.Cells.Item(42, 3).Value = 0
$ws.Cells.Item(42, 4).Value = 10
$ws.Cells.Item(43, 2).Value = "1011100100000000110000000000000010000000000000000001000000000000"
$ws.Cells.Item(43, 3).Value = 0
$ws.Cells.Item(43, 4).Value = 9
$ws.Cells.Item(44, 2).Value = "1011000000110000110000001100000010000000000000000000000000000000"
$ws.Cells.Item(44, 3).Value = 0
$ws.Cells.Item(44, 4).Value = 10
$ws.Cells.Item(45, 2).Value = "1011101100000000110000000000000010000000000000000001000000000000"
$ws.Cells.Item(45, 3).Value = 0.015629
$ws.Cells.Item(45, 4).Value = 8
$ws.Cells.Item(46, 2).Value = "1011100100000000110001001000000000010000000000000000000000000000"
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 10
$ws.Cells.Item(47, 2).Value = "1011100000000000110000001000000000000000000000000000000000000000"
$ws.Cells.Item(47, 3).Value = 0.015615
$ws.Cells.Item(47, 4).Value = 7
$ws.Cells.Item(48, 2).Value = "1011001100000000110000000000000010000000000000000000000000000000"
$ws.Cells.Item(48, 3).Value = 0.005385
$ws.Cells.Item(48, 4).Value = 8
$ws.Cells.Item(49, 2).Value = "1111100100000000100001000001000010000010000000001000000000000000"
$ws.Cells.Item(49, 3).Value = 0.009028
$ws.Cells.Item(49, 4).Value = 8
$ws.Cells.Item(50, 2).Value = "1110100100000000110011000000000000000000000000000010000000000000"
$ws.Cells.Item(50, 3).Value = 0.005836
$ws.Cells.Item(50, 4).Value = 10
$ws.Cells.Item(51, 2).Value = "1011101100000000110000001000000000000000000000000001000000000000"
$ws.Cells.Item(51, 3).Value = 0.011243
$ws.Cells.Item(51, 4).Value = 8
$ws.Cells.Item(52, 2).Value = "1011101100000000110000001000000000000000000000000001000000000000"
$ws.Cells.Item(52, 3).Value = 0.006958
$ws.Cells.Item(52, 4).Value = 8
$ws.Cells.Item(53, 2).Value = "1011100100000000110000000000000000000000000000000001000000000000"
$ws.Cells.Item(53, 3).Value = 0.006878
$ws.Cells.Item(53, 4).Value = 10
$ws.Cells.Item(54, 2).Value = "1010001000000000110000000000000000000000000000000000000000000000"
$ws.Cells.Item(54, 3).Value = 0.013901
$ws.Cells.Item(54, 4).Value = 7
$ws.Cells.Item(55, 2).Value = "1010000000000000100000000000000000010000000000000000000000000000"
$ws.Cells.Item(55, 3).Value = 0.006984
$ws.Cells.Item(55, 4).Value = 10
$ws.Cells.Item(56, 2).Value = "1100111000010000101100000000000000000000100000000000000000000000"
$ws.Cells.Item(56, 3).Value = 0.013855
$ws.Cells.Item(56, 4).Value = 12
$ws.Cells.Item(57, 2).Value = "1011000000000000100000100000000000000000000000001000000000000000"
$ws.Cells.Item(57, 3).Value = 0.01388
$ws.Cells.Item(57, 4).Value = 10
$ws.Cells.Item(58, 2).Value = "1011101100000000110001000000000000010010000000000001000000000000"
$ws.Cells.Item(58, 3).Value = 0.013895
$ws.Cells.Item(58, 4).Value = 10
$ws.Cells.Item(59, 2).Value = "1010000000000000100000100000000010001000000000001000000000000000"
$ws.Cells.Item(59, 3).Value = 0.027684
$ws.Cells.Item(59, 4).Value = 9
$ws.Cells.Item(60, 2).Value = "1010100000010000110000001100000000000010100000011000000000010000"
$ws.Cells.Item(60, 3).Value = 0.020786
$ws.Cells.Item(60, 4).Value = 9
$ws.Cells.Item(61, 2).Value = "1011101001000000111000000000000000000000100000000000010100000000"
$ws.Cells.Item(61, 3).Value = 0.052032
$ws.Cells.Item(61, 4).Value = 12
$ws.Cells.Item(62, 2).Value = "1011101010000000110000000000000011000000000000000000001000010000"
$ws.Cells.Item(62, 3).Value = 0.03027
$ws.Cells.Item(62, 4).Value = 10
$ws.Cells.Item(63, 2).Value = "1011000000000000100000000000000000000000000000000000000000000000"
$ws.Cells.Item(63, 3).Value = 0.046874
$ws.Cells.Item(63, 4).Value = 10
